$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.298.96'
$ws.Range("E2").Value = '  -1.23%  '
$ws.Range("D3").Value = '1.704.31'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").Value = "'223.46"
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").Value = "'0.5306"
$ws.Range("E6").Value = '  -1.27%  '
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").Value = "'0.2653"
$ws.Range("E8").Value = '  -1.71%  '
$ws.Range("D9").Value = "'0.06573"
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("D10").Value = "'20.73"
$ws.Range("E10").Value = '  -4.46%  '
$ws.Range("D11").Value = "'0.07629"
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").Value = "'4.518"
$ws.Range("E12").Value = '  -2.95%  '
$ws.Range("D13").Value = '1.713.49'
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("D14").Value = '1.938.93'
$ws.Range("E14").Value = '  -1.22%  '
$ws.Range("D15").Value = "'0.5748"
$ws.Range("E15").Value = '  -2.34%  '
$ws.Range("D16").Value = '0.0₅8128'
$ws.Range("E16").Value = '  -2.05%  '
$ws.Range("D17").Value = "'67.48"
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("D18").Value = '27.288.67'
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("D19").Value = "'216.50"
$ws.Range("E19").Value = '  -3.36%  '
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("D21").Value = "'4.618"
$ws.Range("E21").Value = '  -2.94%  '
$ws.Range("D22").Value = "'10.39"
$ws.Range("E22").Value = '  -3.26%  '
$ws.Range("D23").Value = "'5.919"
$ws.Range("E23").Value = '  -3.67%  '
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("D25").Value = "'142.57"
$ws.Range("E25").Value = '  -3.79%  '
$ws.Range("D26").Value = "'1.712"
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").Value = "'0.1204"
$ws.Range("E27").Value = '  -2.53%  '
$ws.Range("D28").Value = "'7.206"
$ws.Range("E28").Value = '  -3.05%  '
$ws.Range("D29").Value = "'16.06"
$ws.Range("E29").Value = '  -4.59%  '
$ws.Range("D30").Value = "'0.05367"
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("D31").Value = "'1.286"
$ws.Range("E31").Value = '  -1.36%  '
$ws.Range("D32").Value = "'3.467"
$ws.Range("E32").Value = '  -3.63%  '
$ws.Range("D33").Value = "'3.393"
$ws.Range("E33").Value = '  -2.24%  '
$ws.Range("D34").Value = "'1.637"
$ws.Range("E34").Value = '  -1.67%  '
$ws.Range("D35").Value = "'2.873"
$ws.Range("E35").Value = '  +1.62%  '
$ws.Range("D36").Value = "'2.409"
$ws.Range("E36").Value = '  -1.68%  '
$ws.Range("D37").Value = "'0.9431"
$ws.Range("E37").Value = '  -2.13%  '
$ws.Range("D38").Value = "'0.5810"
$ws.Range("E38").Value = '  -1.88%  '
$ws.Range("D39").Value = "'0.01626"
$ws.Range("E39").Value = '  -1.83%  '
$ws.Range("D40").Value = "'5.749"
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").Value = '1.039.11'
$ws.Range("E42").Value = '  -2.17%  '
$ws.Range("D43").Value = "'0.8384"
$ws.Range("E43").Value = '  -2.49%  '
$ws.Range("D44").Value = "'100.92"
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = '1.847.16'
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("D46").Value = '0.0₈118'
$ws.Range("E46").Value = '  +2.94%  '
$ws.Range("D47").Value = "'57.52"
$ws.Range("E47").Value = '  -2.78%  '
$ws.Range("D48").Value = "'0.4508"
$ws.Range("E48").Value = '  +1.65%  '
$ws.Range("D49").Value = "'1.001"
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").Value = "'8.041"
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("D51").Value = "'0.05228"
$ws.Range("E51").Value = '  -1.15%  '
